# Update cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for the rows that changed, matching the latest scrape.
# For D-column values that are plain decimal numbers (no thousands dots),
# force the cell's number format to Text first so Excel keeps them as
# literal strings (preserving trailing zeros, e.g. "0.190") instead of
# silently re-parsing them as numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.260.97'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.494.93'
$ws.Range('E3').Value = '  -2.58%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.61'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.45'
$ws.Range('E6').Value = '  -2.92%  '
$ws.Range('D7').Value = '3.484.40'
$ws.Range('E7').Value = '  -2.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.611'
$ws.Range('E8').Value = '  -3.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.190'
$ws.Range('E10').Value = '  +4.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.647'
$ws.Range('E11').Value = '  -2.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.98'
$ws.Range('E12').Value = '  -3.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000301'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.42'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').Value = '4.051.57'
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.29'
$ws.Range('E16').Value = '  -3.61%  '
$ws.Range('D17').Value = '69.104.19'
$ws.Range('E17').Value = '  -1.38%  '
$ws.Range('D18').Value = '3.492.18'
$ws.Range('E18').Value = '  -2.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.22'
$ws.Range('E19').Value = '  -3.74%  '
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '540.06'
$ws.Range('E21').Value = '  +12.41%  '
$ws.Range('E22').Value = '  -3.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.30'
$ws.Range('E23').Value = '  -4.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.96'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.42'
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '93.49'
$ws.Range('E26').Value = '  -2.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.23'
$ws.Range('E27').Value = '  +1.07%  '
$ws.Range('E28').Value = '  -2.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.08'
$ws.Range('E29').Value = '  -3.62%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.71'
$ws.Range('E30').Value = '  -1.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.24'
$ws.Range('E31').Value = '  -5.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.58'
$ws.Range('E32').Value = '  +2.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.33'
$ws.Range('E33').Value = '  -3.50%  '
$ws.Range('E34').Value = '  -5.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '537.01'
$ws.Range('E35').Value = '  -8.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '37.81'
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('D40').Value = '0.0₃0762'
$ws.Range('E40').Value = '  -5.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.37'
$ws.Range('E41').Value = '  -2.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.133'
$ws.Range('E42').Value = '  -2.51%  '
$ws.Range('D43').Value = '3.292.20'
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.03'
$ws.Range('E44').Value = '  -7.56%  '
$ws.Range('E45').Value = '  -3.13%  '
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('E47').Value = '  +4.53%  '
$ws.Range('E48').Value = '  -3.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.87'
$ws.Range('E49').Value = '  -6.59%  '
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '137.63'
$ws.Range('E51').Value = '  +2.08%  '
